$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in I1 and J1, matching the style of the existing header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill data rows 2..39: column I is always 1, column J mirrors column H
for ($r = 2; $r -le 39; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
